# [MOSIP-14369] Fix: boolean values
#
# The "is_active" column (D) previously held a TRUE() formula that Excel
# stored as a numeric boolean (t="n" / t="b" with <f>TRUE()</f>). The fix
# replaces that with the literal text string "TRUE" (shared string, t="s"),
# keeping the existing text-formatted cell style.
#
# Assigning the bare string "TRUE" via Range.Value/.Value2 gets auto-coerced
# by this engine into a real boolean cell (t="b"), and prefixing with an
# apostrophe forces a "quote-prefixed" text style (a new, different cell
# style) instead of reusing the original style. To land the literal text
# "TRUE" while preserving the original style (s="1") exactly, we stage the
# text in a scratch cell as a formula that evaluates to the text "TRUE",
# copy it, and paste-special *values only* into the target cells - which
# carries over the text value without touching the destination formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# Scratch cell, far away from the used range, used only to mint a literal
# text "TRUE" value we can copy from.
$helper = $ws.Range("Z1")
$helper.Formula = "=""TRUE"""
$helper.Copy() | Out-Null

$ws.Range("D2").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("D3").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("D4").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("D5").PasteSpecial($xlPasteValues) | Out-Null

$helper.ClearContents()
$excel.CutCopyMode = 0

# Matches the updated <selection> in the saved sheet view.
$ws.Range("I19").Select() | Out-Null
